$wb = $excel.ActiveWorkbook

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H62").Value2 = 7062.857
$ws.Range("I62").Value2 = 6903
$ws.Range("K62").Value2 = 6903
$ws.Range("M62").Value2 = -6279
$ws.Range("H65").Value2 = 7062.857
$ws.Range("I65").Value2 = 6903
$ws.Range("K65").Value2 = 34515
$ws.Range("M65").Value2 = -31395
$ws.Range("H76").Value2 = 5490.857
$ws.Range("J76").Value2 = 7999.5
$ws.Range("L76").Value2 = 7999.5
$ws.Range("N76").Value2 = -8629.5
$ws.Range("H79").Value2 = 5490.857
$ws.Range("J79").Value2 = 7999.5
$ws.Range("L79").Value2 = 7999.5
$ws.Range("N79").Value2 = -10183.5
$ws.Range("H86").Value2 = 250012350
$ws.Range("J86").Value2 = 2399.5
$ws.Range("L86").Value2 = 2399.5
$ws.Range("N86").Value2 = -4645.5
$ws.Range("H89").Value2 = 250012350
$ws.Range("J89").Value2 = 2399.5
$ws.Range("L89").Value2 = 11997.5
$ws.Range("N89").Value2 = -23229.5
$ws.Range("H98").Value2 = 3050.025
$ws.Range("J98").Value2 = 1499.6
$ws.Range("L98").Value2 = 1499.6
$ws.Range("N98").Value2 = -4495.6
$ws.Range("H107").Value2 = 564.73334
$ws.Range("I107").Value2 = 564.73334
$ws.Range("K107").Value2 = 564.73334
$ws.Range("M107").Value2 = 1355.26666
$ws.Range("H118").Value2 = 946.3125
$ws.Range("I118").Value2 = 724.4286
$ws.Range("K118").Value2 = 2173.2858
$ws.Range("M118").Value2 = -516.2857999999997
$ws.Range("H122").Value2 = 3050.025
$ws.Range("J122").Value2 = 1499.6
$ws.Range("L122").Value2 = 4498.799999999999
$ws.Range("N122").Value2 = -9398.799999999999
$ws.Range("H126").Value2 = 77751.336
$ws.Range("J126").Value2 = 77751.336
$ws.Range("L126").Value2 = 77751.336
$ws.Range("N126").Value2 = -87631.336
$ws.Range("H132").Value2 = 4491.5117
$ws.Range("I132").Value2 = 4574.643
$ws.Range("K132").Value2 = 13723.929
$ws.Range("M132").Value2 = -11193.929
$ws.Range("H135").Value2 = 547.76666
$ws.Range("I135").Value2 = 549.4074000000001
$ws.Range("K135").Value2 = 4944.6666
$ws.Range("M135").Value2 = -2409.6666
$ws.Range("H137").Value2 = 3576.6667
$ws.Range("I137").Value2 = 3228.6667
$ws.Range("J137").Value2 = 4272.6665
$ws.Range("K137").Value2 = 9686.000100000001
$ws.Range("L137").Value2 = 12817.9995
$ws.Range("M137").Value2 = -7136.000100000001
$ws.Range("N137").Value2 = -17917.9995
$ws.Range("H138").Value2 = 454032.3
$ws.Range("I138").Value2 = 5570.4287
$ws.Range("J138").Value2 = 518098.28
$ws.Range("K138").Value2 = 16711.2861
$ws.Range("L138").Value2 = 1554294.84
$ws.Range("M138").Value2 = -11571.2861
$ws.Range("N138").Value2 = -1564574.84
$ws.Range("H141").Value2 = 3734.5386
$ws.Range("I141").Value2 = 3492.5715
$ws.Range("J141").Value2 = 4016.8333
$ws.Range("K141").Value2 = 10477.7145
$ws.Range("L141").Value2 = 12050.4999
$ws.Range("M141").Value2 = -5297.7145
$ws.Range("N141").Value2 = -22410.4999

# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H3").Value2 = 377873.5
$ws.Range("J3").Value2 = 11500
$ws.Range("L3").Value2 = 11500
$ws.Range("N3").Value2 = -11730
$ws.Range("H32").Value2 = 37170.832
$ws.Range("I32").Value2 = 11143.862
$ws.Range("K32").Value2 = 11143.862
$ws.Range("M32").Value2 = -10856.862
$ws.Range("H61").Value2 = 6540.8237
$ws.Range("I61").Value2 = 2686.5557
$ws.Range("J61").Value2 = 10876.875
$ws.Range("K61").Value2 = 2686.5557
$ws.Range("L61").Value2 = 10876.875
$ws.Range("M61").Value2 = -2474.5557
$ws.Range("N61").Value2 = -11300.875
$ws.Range("H74").Value2 = 141320.62
$ws.Range("I74").Value2 = 187062.1
$ws.Range("J74").Value2 = 4096.2
$ws.Range("K74").Value2 = 187062.1
$ws.Range("L74").Value2 = 4096.2
$ws.Range("M74").Value2 = -186188.1
$ws.Range("N74").Value2 = -5844.2
$ws.Range("H77").Value2 = 141320.62
$ws.Range("I77").Value2 = 187062.1
$ws.Range("J77").Value2 = 4096.2
$ws.Range("K77").Value2 = 935310.5
$ws.Range("L77").Value2 = 20481
$ws.Range("M77").Value2 = -930942.5
$ws.Range("N77").Value2 = -29217
$ws.Range("H132").Value2 = 2028.0605
$ws.Range("I132").Value2 = 1642.7742
$ws.Range("J132").Value2 = 8000
$ws.Range("K132").Value2 = 4928.3226
$ws.Range("L132").Value2 = 24000
$ws.Range("M132").Value2 = -2398.3226
$ws.Range("N132").Value2 = -29060
$ws.Range("H136").Value2 = 6540.8237
$ws.Range("I136").Value2 = 2686.5557
$ws.Range("J136").Value2 = 10876.875
$ws.Range("K136").Value2 = 8059.6671
$ws.Range("L136").Value2 = 32630.625
$ws.Range("M136").Value2 = -5509.6671
$ws.Range("N136").Value2 = -37730.625
$ws.Range("H139").Value2 = 91626.14
$ws.Range("J139").Value2 = 91626.14
$ws.Range("L139").Value2 = 91626.14
$ws.Range("N139").Value2 = -101906.14

# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H82").Value2 = 51140.5
$ws.Range("I82").Value2 = 24998
$ws.Range("J82").Value2 = 77283
$ws.Range("K82").Value2 = 24998
$ws.Range("L82").Value2 = 77283
$ws.Range("M82").Value2 = -24615
$ws.Range("N82").Value2 = -78049
$ws.Range("H85").Value2 = 51140.5
$ws.Range("I85").Value2 = 24998
$ws.Range("J85").Value2 = 77283
$ws.Range("K85").Value2 = 24998
$ws.Range("L85").Value2 = 77283
$ws.Range("M85").Value2 = -23672
$ws.Range("N85").Value2 = -79935
$ws.Range("H86").Value2 = 3425
$ws.Range("I86").Value2 = 2741.7778
$ws.Range("K86").Value2 = 2741.7778
$ws.Range("M86").Value2 = -1618.7778
$ws.Range("H89").Value2 = 3425
$ws.Range("I89").Value2 = 2741.7778
$ws.Range("K89").Value2 = 13708.889
$ws.Range("M89").Value2 = -8092.888999999999
$ws.Range("H107").Value2 = 1716.8108
$ws.Range("I107").Value2 = 1387.2727
$ws.Range("J107").Value2 = 2200.1333
$ws.Range("K107").Value2 = 1387.2727
$ws.Range("L107").Value2 = 2200.1333
$ws.Range("M107").Value2 = 532.7273
$ws.Range("N107").Value2 = -6040.1333
$ws.Range("H134").Value2 = 2731.3809
$ws.Range("I134").Value2 = 2535
$ws.Range("J134").Value2 = 3222.3333
$ws.Range("K134").Value2 = 7605
$ws.Range("L134").Value2 = 9666.999899999999
$ws.Range("M134").Value2 = -5070
$ws.Range("N134").Value2 = -14736.9999

# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value2 = 1998.3334
$ws.Range("I22").Value2 = 1995.3334
$ws.Range("J22").Value2 = 1999.8334
$ws.Range("K22").Value2 = 1995.3334
$ws.Range("L22").Value2 = 1999.8334
$ws.Range("M22").Value2 = -1645.3334
$ws.Range("N22").Value2 = -2699.8334
$ws.Range("H31").Value2 = 3705.3215
$ws.Range("I31").Value2 = 2997.9
$ws.Range("J31").Value2 = 4098.3335
$ws.Range("K31").Value2 = 2997.9
$ws.Range("L31").Value2 = 4098.3335
$ws.Range("M31").Value2 = -2702.9
$ws.Range("N31").Value2 = -4688.3335
$ws.Range("H34").Value2 = 3705.3215
$ws.Range("I34").Value2 = 2997.9
$ws.Range("J34").Value2 = 4098.3335
$ws.Range("K34").Value2 = 2997.9
$ws.Range("L34").Value2 = 4098.3335
$ws.Range("M34").Value2 = -2795.9
$ws.Range("N34").Value2 = -4502.3335
$ws.Range("H58").Value2 = 4109.5
$ws.Range("J58").Value2 = 4714.2856
$ws.Range("L58").Value2 = 4714.2856
$ws.Range("N58").Value2 = -5120.2856
$ws.Range("H99").Value2 = 5529.273
$ws.Range("I99").Value2 = 5690.375
$ws.Range("K99").Value2 = 5690.375
$ws.Range("M99").Value2 = -4192.375
$ws.Range("H105").Value2 = 1790.6875
$ws.Range("I105").Value2 = 1089.2273
$ws.Range("J105").Value2 = 3333.9
$ws.Range("K105").Value2 = 1089.2273
$ws.Range("L105").Value2 = 3333.9
$ws.Range("M105").Value2 = 657.7727
$ws.Range("N105").Value2 = -6827.9
$ws.Range("H107").Value2 = 549.2
$ws.Range("I107").Value2 = 499.1111
$ws.Range("K107").Value2 = 499.1111
$ws.Range("M107").Value2 = 1420.8889
$ws.Range("H126").Value2 = 5529.273
$ws.Range("I126").Value2 = 5690.375
$ws.Range("K126").Value2 = 17071.125
$ws.Range("M126").Value2 = -14601.125
$ws.Range("H132").Value2 = 4362.625
$ws.Range("I132").Value2 = 4172.343
$ws.Range("J132").Value2 = 4874.923
$ws.Range("K132").Value2 = 12517.029
$ws.Range("L132").Value2 = 14624.769
$ws.Range("M132").Value2 = -9987.028999999999
$ws.Range("N132").Value2 = -19684.769
$ws.Range("H134").Value2 = 3052.2
$ws.Range("I134").Value2 = 2177.6
$ws.Range("K134").Value2 = 6532.799999999999
$ws.Range("M134").Value2 = -3997.799999999999
$ws.Range("H136").Value2 = 4109.5
$ws.Range("J136").Value2 = 4714.2856
$ws.Range("L136").Value2 = 14142.8568
$ws.Range("N136").Value2 = -19242.8568
$ws.Range("H141").Value2 = 1896040.8
$ws.Range("J141").Value2 = 1896040.8
$ws.Range("L141").Value2 = 1896040.8
$ws.Range("N141").Value2 = -1906400.8

# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value2 = 1851.6154
$ws.Range("I5").Value2 = 567.5185
$ws.Range("K5").Value2 = 1702.5555
$ws.Range("M5").Value2 = -1590.5555
$ws.Range("H80").Value2 = 3724.75
$ws.Range("J80").Value2 = 4299.6665
$ws.Range("L80").Value2 = 12898.9995
$ws.Range("N80").Value2 = -14770.9995
$ws.Range("H83").Value2 = 3724.75
$ws.Range("J83").Value2 = 4299.6665
$ws.Range("L83").Value2 = 38696.9985
$ws.Range("N83").Value2 = -48056.9985
$ws.Range("H98").Value2 = 448.57144
$ws.Range("J98").Value2 = 349.66666
$ws.Range("L98").Value2 = 1048.99998
$ws.Range("N98").Value2 = -4044.99998
$ws.Range("H107").Value2 = 1248
$ws.Range("I107").Value2 = 1084.3334
$ws.Range("J107").Value2 = 1283.0714
$ws.Range("K107").Value2 = 3253.0002
$ws.Range("L107").Value2 = 3849.2142
$ws.Range("M107").Value2 = -1333.0002
$ws.Range("N107").Value2 = -7689.2142
$ws.Range("H113").Value2 = 4330.6924
$ws.Range("J113").Value2 = 4458.25
$ws.Range("L113").Value2 = 13374.75
$ws.Range("N113").Value2 = -17714.75
$ws.Range("H122").Value2 = 2201.1428
$ws.Range("I122").Value2 = 1825.6666
$ws.Range("J122").Value2 = 2482.75
$ws.Range("K122").Value2 = 16430.9994
$ws.Range("L122").Value2 = 22344.75
$ws.Range("M122").Value2 = -13980.9994
$ws.Range("N122").Value2 = -27244.75
$ws.Range("H135").Value2 = 1851.6154
$ws.Range("I135").Value2 = 567.5185
$ws.Range("K135").Value2 = 5107.6665
$ws.Range("M135").Value2 = -2572.6665

# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H18").Value2 = 0
$ws.Range("J18").Value2 = 0
$ws.Range("L18").Value2 = 0
$ws.Range("N18").ClearContents()
$ws.Range("H46").Value2 = 17343.666
$ws.Range("I46").Value2 = 17343.666
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 17343.666
$ws.Range("L46").Value2 = 0
$ws.Range("M46").Value2 = -17187.666
$ws.Range("N46").ClearContents()
$ws.Range("I80").Value2 = 200001500
$ws.Range("J80").Value2 = 2666.5
$ws.Range("K80").Value2 = 200001500
$ws.Range("L80").Value2 = 2666.5
$ws.Range("M80").Value2 = -200000502
$ws.Range("N80").Value2 = -4662.5
$ws.Range("I83").Value2 = 200001500
$ws.Range("J83").Value2 = 2666.5
$ws.Range("K83").Value2 = 1000007500
$ws.Range("L83").Value2 = 13332.5
$ws.Range("M83").Value2 = -1000002508
$ws.Range("N83").Value2 = -23316.5
$ws.Range("H97").Value2 = 789.7826
$ws.Range("I97").Value2 = 788.5294
$ws.Range("J97").Value2 = 793.3333
$ws.Range("K97").Value2 = 788.5294
$ws.Range("L97").Value2 = 793.3333
$ws.Range("M97").Value2 = -292.5294
$ws.Range("N97").Value2 = -1785.3333
$ws.Range("H107").Value2 = 394.29413
$ws.Range("J107").Value2 = 423.66666
$ws.Range("L107").Value2 = 423.66666
$ws.Range("N107").Value2 = -4263.66666
$ws.Range("H126").Value2 = 10958.8
$ws.Range("I126").Value2 = 4727.0713
$ws.Range("K126").Value2 = 14181.2139
$ws.Range("M126").Value2 = -11711.2139
$ws.Range("H132").Value2 = 2572.1
$ws.Range("I132").Value2 = 1694.2609
$ws.Range("K132").Value2 = 5082.7827
$ws.Range("M132").Value2 = -2552.7827

# --- Sheet 7 (LTW) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value2 = 41902.97
$ws.Range("I40").Value2 = 207798.6
$ws.Range("K40").Value2 = 207798.6
$ws.Range("M40").Value2 = -207662.6
$ws.Range("H56").Value2 = 0
$ws.Range("I56").Value2 = 0
$ws.Range("K56").Value2 = 0
$ws.Range("M56").ClearContents()
$ws.Range("J99").Value2 = 25000
$ws.Range("L99").Value2 = 25000
$ws.Range("N99").Value2 = -30990
$ws.Range("H136").Value2 = 6714.2856
$ws.Range("I136").Value2 = 3500
$ws.Range("J136").Value2 = 11000
$ws.Range("K136").Value2 = 10500
$ws.Range("L136").Value2 = 33000
$ws.Range("M136").Value2 = -7950
$ws.Range("N136").Value2 = -38100

# --- Sheet 8 (WVR) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H51").Value2 = 34000
$ws.Range("J51").Value2 = 35000
$ws.Range("L51").Value2 = 35000
$ws.Range("N51").Value2 = -36020
$ws.Range("H62").Value2 = 8999.5
$ws.Range("J62").Value2 = 8999.5
$ws.Range("L62").Value2 = 8999.5
$ws.Range("N62").Value2 = -10247.5
$ws.Range("H65").Value2 = 8999.5
$ws.Range("J65").Value2 = 8999.5
$ws.Range("L65").Value2 = 44997.5
$ws.Range("N65").Value2 = -51237.5
$ws.Range("H81").Value2 = 8048.8335
$ws.Range("I81").Value2 = 5658.6
$ws.Range("K81").Value2 = 11317.2
$ws.Range("M81").Value2 = -10256.2
$ws.Range("H84").Value2 = 8048.8335
$ws.Range("I84").Value2 = 5658.6
$ws.Range("K84").Value2 = 56586
$ws.Range("M84").Value2 = -51282
$ws.Range("H132").Value2 = 6533.1387
$ws.Range("I132").Value2 = 6173.1
$ws.Range("K132").Value2 = 18519.3
$ws.Range("M132").Value2 = -15989.3
$ws.Range("H136").Value2 = 2221.6428
$ws.Range("I136").Value2 = 1693.5
$ws.Range("J136").Value2 = 2925.8333
$ws.Range("K136").Value2 = 5080.5
$ws.Range("L136").Value2 = 8777.499899999999
$ws.Range("M136").Value2 = -2530.5
$ws.Range("N136").Value2 = -13877.4999
